$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33
$srcRow = 32

# Text columns (A-D): force text interpretation (so date/time-like strings
# and the zero-padded week number aren't auto-converted to numbers/dates),
# then copy the formatting (style) from the row above so no explicit style
# index ends up on the new cells, matching the existing rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-08"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "14:36:27"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "Monday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "01"

$ws.Range("A$srcRow`:D$srcRow").Copy()
$ws.Range("A$row`:D$row").PasteSpecial(-4122)

# Numeric columns (E-T)
$ws.Cells.Item($row, 5).Value = 139556
$ws.Cells.Item($row, 6).Value = 143110
$ws.Cells.Item($row, 7).Value = 172278
$ws.Cells.Item($row, 8).Value = 147229
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118267
$ws.Cells.Item($row, 11).Value = 224722
$ws.Cells.Item($row, 12).Value = 249715
$ws.Cells.Item($row, 13).Value = 185267
$ws.Cells.Item($row, 14).Value = 110479
$ws.Cells.Item($row, 15).Value = 40644
$ws.Cells.Item($row, 16).Value = 30815
$ws.Cells.Item($row, 17).Value = 72412
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42094
$ws.Cells.Item($row, 20).Value = -1
